$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1) - clone formatting from the existing
# header cell (G1) so it reuses the same bold/bordered header style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value for the Save column (H2), left unstyled like the other
# numeric data cells in row 2.
$ws.Range("H2").Value = 0
